$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.700.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.525.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.549.90'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.364'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.972.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.36%  '
$ws.Range("E15").Value = '  -2.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.522.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.535.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  -2.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.435'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.26%  '
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.17%  '
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("E32").Value = '  -5.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.23%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.41'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.96%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '293.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.603'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0515'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0228'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.60%  '
